$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Hüllflächen, Himmelsricht."  (zone envelope areas)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 5 - Zone 1 North: drop the "=12.5*10" formula, typed value instead
$ws1.Range("Y5").Value = 136.30000000000001

# Row 6 - Zone 1 East: formula changed to match the one used in V5
# (copy V5's format first so the border matches, then overwrite the formula)
$ws1.Range("V5").Copy()
$ws1.Range("V6").PasteSpecial(-4122)
$ws1.Range("V6").Formula = "=25*16.66"
$ws1.Range("Y6").Value = 126.93

# Row 8 - Zone 1 West: both formulas replaced by typed values
$ws1.Range("V8").Value = 146.52000000000001
$ws1.Range("Y8").Value = 37.146000000000001

# Row 9 - Zone 2 North: formula replaced by typed value
$ws1.Range("V9").Value = 43.316000000000003

# Row 10 - Zone 2 East: both formulas replaced by typed values
$ws1.Range("V10").Value = 942.95600000000002
$ws1.Range("Y10").Value = 126.8

# Row 11 - Zone 2 South: typed value updated
$ws1.Range("V11").Value = 301.14999999999998

# ---------------------------------------------------------------------------
# Sheet 2: "Strukturen Hüllfläche" (wall/structure reference tables)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# The swimming-pool wall thickness entry moves from 16 cm to 1600 (mm) -
# the previous reference row is preserved off to the side (columns Q:U).
$ws2.Range("E4:I4").Copy($ws2.Range("Q4:U4"))
$ws2.Range("E5:I5").Copy($ws2.Range("Q5:U5"))
$ws2.Range("E6:I6").Copy($ws2.Range("Q6:U6"))
$ws2.Range("E7:I7").Copy($ws2.Range("Q7:U7"))

$ws2.Range("E5").Value = 1600

$ws2.Range("E22:I22").Copy($ws2.Range("Q22:U22"))
$ws2.Range("E23:I23").Copy($ws2.Range("Q23:U23"))
$ws2.Range("E24:I24").Copy($ws2.Range("Q24:U24"))
$ws2.Range("E25:I25").Copy($ws2.Range("Q25:U25"))
$ws2.Range("E26:I26").Copy($ws2.Range("Q26:U26"))
$ws2.Range("E27:I27").Copy($ws2.Range("Q27:U27"))
